$wb = $excel.ActiveWorkbook

# --- Logs sheet: append new row 11 with the latest mail log entry ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(11, 1).Value = "Retour status"
$logs.Cells.Item(11, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(11, 4).Value = "Retour / Terugbetaling"
$logs.Cells.Item(11, 6).Value = "2025-08-26 21:18:17"
$logs.Cells.Item(11, 7).Value = "Nee"
$logs.Cells.Item(11, 8).Value = "Ja"
$logs.Cells.Item(11, 9).Value = "Nee"
$logs.Cells.Item(11, 10).Value = "Nee"

# Extend the conditional formatting ranges to include the new row
$logs.Range("D2:D10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D11"))
$logs.Range("G2:G10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G11"))
$logs.Range("H2:H10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H11"))
$logs.Range("I2:I10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I11"))
$logs.Range("J2:J10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J11"))

# --- Dashboard sheet: bump the "Retour / Terugbetaling" count ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(2, 2).Value = 6
